$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Dashboard")
$ws.Activate()
